# Add coefficient tables (d_name / mu_name lookups) from multistate NMA and MA
# to each of the 4 parameter-lookup sheets.

$wb = $excel.ActiveWorkbook

$weibull    = $wb.Worksheets.Item("weibull")
$gompertz   = $wb.Worksheets.Item("gompertz")
$fracpoly1  = $wb.Worksheets.Item("fracpoly1")
$fracpoly2  = $wb.Worksheets.Item("fracpoly2")

# --- weibull (sheet1) ---------------------------------------------------
$weibull.Range("G1").Value = "d_name"
$weibull.Range("H1").Value = "mu_name"
$weibull.Range("G2").Value = "d_1"
$weibull.Range("G3").Value = "d_2"
$weibull.Range("H2").Value = "mu_1"
$weibull.Range("H3").Value = "mu_2"

# mu_3 is first introduced on fracpoly1 so the shared-string table ends up
# in the same order as the authored workbook.
$fracpoly1.Range("H4").Value = "mu_3"

$weibull.Range("H6").Value = "mu_5"
$weibull.Range("H7").Value = "mu_6"
$weibull.Range("H4").Value = "mu_4"

# --- gompertz (sheet2) ---------------------------------------------------
$gompertz.Range("G1").Value = "d_name"
$gompertz.Range("H1").Value = "mu_name"
$gompertz.Range("G2").Value = "d_1"
$gompertz.Range("H2").Value = "mu_1"
$gompertz.Range("G3").Value = "d_2"
$gompertz.Range("H3").Value = "mu_2"
$gompertz.Range("H4").Value = "mu_4"
$gompertz.Range("H6").Value = "mu_5"
$gompertz.Range("H7").Value = "mu_6"

# --- fracpoly1 (sheet3) ---------------------------------------------------
$fracpoly1.Range("G1").Value = "d_name"
$fracpoly1.Range("H1").Value = "mu_name"
$fracpoly1.Range("G2").Value = "d_1"
$fracpoly1.Range("H2").Value = "mu_1"
$fracpoly1.Range("G3").Value = "d_2"
$fracpoly1.Range("H3").Value = "mu_2"
$fracpoly1.Range("H5").Value = "mu_4"
$fracpoly1.Range("H8").Value = "mu_5"
$fracpoly1.Range("H9").Value = "mu_6"

# --- fracpoly2 (sheet4) ---------------------------------------------------
$fracpoly2.Range("G1").Value = "d_name"
$fracpoly2.Range("H1").Value = "mu_name"
$fracpoly2.Range("G2").Value = "d_1"
$fracpoly2.Range("H2").Value = "mu_1"
$fracpoly2.Range("G3").Value = "d_2"
$fracpoly2.Range("H3").Value = "mu_2"
$fracpoly2.Range("H4").Value = "mu_3"
$fracpoly2.Range("H5").Value = "mu_4"
$fracpoly2.Range("H8").Value = "mu_5"
$fracpoly2.Range("H9").Value = "mu_6"

# --- selections / active sheet -------------------------------------------
# Final selection state on each sheet (order matters: the last Select()
# call below determines which sheet stays active/tabSelected).
$weibull.Range("H6").Select()
$fracpoly1.Range("G1:H10").Select()
$fracpoly2.Range("I8").Select()
$gompertz.Range("H7").Select()
